$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain plain text so numeric-looking price
# strings such as 0.9964 are not auto-converted to numbers by Excel.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '24.011.58'
$ws.Range('E2').Value = '  -1.69%  '
$ws.Range('D3').Value = '1.627.30'
$ws.Range('E3').Value = '  -1.67%  '
$ws.Range('D4').Value = '0.9964'
$ws.Range('E4').Value = '  -0.43%  '
$ws.Range('D5').Value = '309.78'
$ws.Range('E5').Value = '  -0.59%  '
$ws.Range('D6').Value = '0.9973'
$ws.Range('E6').Value = '  -0.39%  '
$ws.Range('D7').Value = '0.3935'
$ws.Range('E7').Value = '  +0.46%  '
$ws.Range('D8').Value = '0.3858'
$ws.Range('E8').Value = '  -1.19%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = '50.28'
$ws.Range('E9').Value = '  +0.43%  '
$ws.Range('B10').Value = 'BinanceUSD'
$ws.Range('C10').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D10').Value = '0.9965'
$ws.Range('E10').Value = '  -0.45%  '
$ws.Range('D11').Value = '1.368'
$ws.Range('E11').Value = '  -1.03%  '
$ws.Range('D12').Value = '0.08502'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('D13').Value = '24.03'
$ws.Range('E13').Value = '  -3.57%  '
$ws.Range('D14').Value = '7.114'
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('D15').Value = '7.654'
$ws.Range('E15').Value = '  +0.62%  '
$ws.Range('D16').Value = '0.00001291'
$ws.Range('E16').Value = '  -1.08%  '
$ws.Range('D17').Value = '1.626.54'
$ws.Range('E17').Value = '  -2.30%  '
$ws.Range('D18').Value = '94.12'
$ws.Range('E18').Value = '  +0.99%  '
$ws.Range('D19').Value = '0.06929'
$ws.Range('E19').Value = '  -0.31%  '
$ws.Range('D20').Value = '20.20'
$ws.Range('E20').Value = '  -3.55%  '
$ws.Range('D21').Value = '6.883'
$ws.Range('E21').Value = '  -1.70%  '
$ws.Range('D22').Value = '0.9973'
$ws.Range('E22').Value = '  -0.42%  '
$ws.Range('D23').Value = '13.51'
$ws.Range('E23').Value = '  -2.11%  '
$ws.Range('D24').Value = '24.018.62'
$ws.Range('E24').Value = '  -1.62%  '
$ws.Range('D25').Value = '2.462'
$ws.Range('E25').Value = '  +5.45%  '
$ws.Range('D26').Value = '2.894'
$ws.Range('E26').Value = '  +3.82%  '
$ws.Range('D27').Value = '22.31'
$ws.Range('E27').Value = '  -1.63%  '
$ws.Range('D28').Value = '157.05'
$ws.Range('E28').Value = '  -1.09%  '
$ws.Range('D29').Value = '140.93'
$ws.Range('D30').Value = '5.311'
$ws.Range('E30').Value = '  -7.27%  '
$ws.Range('D31').Value = '7.975'
$ws.Range('E31').Value = '  -2.56%  '
$ws.Range('D32').Value = '2.490'
$ws.Range('E32').Value = '  -1.63%  '
$ws.Range('D33').Value = '1.800.48'
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('D34').Value = '0.08198'
$ws.Range('E34').Value = '  +1.14%  '
$ws.Range('D35').Value = '0.9991'
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('D36').Value = '0.02920'
$ws.Range('D37').Value = '6.676'
$ws.Range('E37').Value = '  -2.57%  '
$ws.Range('D38').Value = '0.2684'
$ws.Range('E38').Value = '  -2.86%  '
$ws.Range('D39').Value = '10.57'
$ws.Range('E39').Value = '  +3.64%  '
$ws.Range('D40').Value = '0.09178'
$ws.Range('E40').Value = '  -2.81%  '
$ws.Range('D41').Value = '13.73'
$ws.Range('E41').Value = '  +2.85%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.7599'
$ws.Range('E42').Value = '  -2.32%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '1.435'
$ws.Range('E43').Value = '  -3.49%  '
$ws.Range('D44').Value = '16.07'
$ws.Range('E44').Value = '  -1.65%  '
$ws.Range('D45').Value = '0.6978'
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('D46').Value = '2.486'
$ws.Range('E46').Value = '  -2.57%  '
$ws.Range('D47').Value = '4.089'
$ws.Range('E47').Value = '  -1.39%  '
$ws.Range('D48').Value = '0.9973'
$ws.Range('E48').Value = '  -0.35%  '
$ws.Range('D49').Value = '0.08314'
$ws.Range('E49').Value = '  -3.10%  '
$ws.Range('D50').Value = '136.73'
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('D51').Value = '1.216'
$ws.Range('E51').Value = '  -6.67%  '
